# Add week 47 report (week 46 data)
# Populates row 21 (week 46) across the Hospital / ICU / Hospital_Regions /
# ICU_Regions data sheets, and lightly revises a handful of previously
# reported values (weeks 43-45) that were updated in the source report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Hospital sheet
# ---------------------------------------------------------------------
$wsHospital = $wb.Worksheets.Item("Hospital")
$wsHospital.Range("B18").Value = 10.44
$wsHospital.Range("B19").Value = 13.54
$wsHospital.Range("B20").Value = 14.23
$wsHospital.Range("B21").Value = 16.739999999999998
$wsHospital.Range("C21").Value = 0

# ---------------------------------------------------------------------
# ICU sheet
# ---------------------------------------------------------------------
$wsICU = $wb.Worksheets.Item("ICU")
$wsICU.Range("B20").Value = 1
$wsICU.Range("B21").Value = 1.05
$wsICU.Range("C21").Value = 0

# ---------------------------------------------------------------------
# Hospital_Regions sheet
# ---------------------------------------------------------------------
$wsHospReg = $wb.Worksheets.Item("Hospital_Regions")
$wsHospReg.Range("D18").Value = 6
$wsHospReg.Range("D19").Value = 8.34
$wsHospReg.Range("I19").Value = 18.170000000000002
$wsHospReg.Range("C20").Value = 7.79
$wsHospReg.Range("D20").Value = 7.02
$wsHospReg.Range("I20").Value = 17.71
$wsHospReg.Range("B21").Value = 24.64
$wsHospReg.Range("C21").Value = 9.7100000000000009
$wsHospReg.Range("D21").Value = 7.78
$wsHospReg.Range("E21").Value = 26.1
$wsHospReg.Range("F21").Value = 20.64
$wsHospReg.Range("G21").Value = 8.9600000000000009
$wsHospReg.Range("H21").Value = 15.47
$wsHospReg.Range("I21").Value = 21.05
$wsHospReg.Range("J21").Value = 29.41

# ---------------------------------------------------------------------
# ICU_Regions sheet
# ---------------------------------------------------------------------
$wsICUReg = $wb.Worksheets.Item("ICU_Regions")
$wsICUReg.Range("F19").Value = 0.95
$wsICUReg.Range("D20").Value = 0.92
$wsICUReg.Range("F20").Value = 0.89
$wsICUReg.Range("I20").Value = 1.79
$wsICUReg.Range("B21").Value = 1.73
$wsICUReg.Range("C21").Value = 1.75
$wsICUReg.Range("D21").Value = 0.63
$wsICUReg.Range("E21").Value = 1.26
$wsICUReg.Range("F21").Value = 1.32
$wsICUReg.Range("G21").Value = 0.65
$wsICUReg.Range("H21").Value = 0.66
$wsICUReg.Range("I21").Value = 1.52
$wsICUReg.Range("J21").Value = 1.45

# ---------------------------------------------------------------------
# View / selection state.
# Selecting a range activates its parent sheet, so the sheets that must
# NOT end up as the active tab are handled first; "Hospital" (the sheet
# that should end active) is selected/activated last.
# ---------------------------------------------------------------------
$wsICU.Range("C15:C21").Select()

$wsHospReg.Range("B2").Select()

$wsICUReg.Range("B2").Select()

$wsICURegGraph = $wb.Worksheets.Item("ICU_Regions_Graph")
$wsICURegGraph.Range("A40").Select()

$wsHospital.Range("B22").Select()
$wsHospital.Activate()
